$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new hourly data point for 2026/01/07 was collected (commit: "daily auto
# push: 2026-01-07 06:50 UTC"). It is inserted as a brand-new row right
# before the 2026/12/29 block (current row 573), which pushes every
# subsequent row down by one (old 573..614 -> new 574..615) and grows the
# used range from D614 to D615.
$ws.Rows(573).Insert()

# Column A holds dates as plain text (e.g. "2026/01/07"), not real Excel
# date serials. Force text formatting before assigning so the COM layer
# doesn't auto-coerce the "yyyy/mm/dd"-looking string into a date value,
# then drop back to the default style so no stray formatting is left on
# the cell (matching the unstyled neighbouring date cells).
$ws.Range("A573").NumberFormat = "@"
$ws.Range("A573").Value = "2026/01/07"
$ws.Range("A573").Style = "Normal"

$ws.Range("B573").Value = "水"
$ws.Range("C573").Value = 14
$ws.Range("D573").Value = 23
